# Update the "Förändrad" (Changed) date column (column C) from 45178 (2023-09-09)
# to 45179 (2023-09-10) for every data row in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
